$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new column M for year 2022, matching the
# formatting of the existing column L (copy format, then set the value).

# Row 3: empty bottom-border cell
$ws.Cells.Item(3, 12).Copy($ws.Cells.Item(3, 13))

# Row 4: header year value
$ws.Cells.Item(4, 12).Copy($ws.Cells.Item(4, 13))
$ws.Cells.Item(4, 13).Value = 2022

# Row 6: Mammals value
$ws.Cells.Item(6, 12).Copy($ws.Cells.Item(6, 13))
$ws.Cells.Item(6, 13).Value = 18

# Row 7: Birds value
$ws.Cells.Item(7, 12).Copy($ws.Cells.Item(7, 13))
$ws.Cells.Item(7, 13).Value = 6.2

# Row 8: Amphibians and Reptiles value (same "-" placeholder used elsewhere in the row)
$ws.Cells.Item(8, 12).Copy($ws.Cells.Item(8, 13))
$ws.Cells.Item(8, 13).Value = "-"

# Move the active selection to match the saved worksheet state
$ws.Range("N4").Select()
